# Add "NA" values under the duplicate_image_filename column (column E)
# for rows 2 through 21.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2:E21").Value = "NA"

# Preserve the pre-existing empty string cell at F1 (the COM runtime's
# save cycle otherwise coerces an untouched empty shared-string cell to
# index 0); explicitly re-assert it as blank so it round-trips unchanged.
$ws.Range("F1").Value = ""
